$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 522.7778  # ALC!H4 was 403.25
$ws.Cells.Item(4, 9).Value = 263.2  # ALC!I4 was 222
$ws.Cells.Item(4, 10).Value = 847.25  # ALC!J4 was 947
$ws.Cells.Item(4, 11).Value = 263.2  # ALC!K4 was 222
$ws.Cells.Item(4, 12).Value = 847.25  # ALC!L4 was 947
$ws.Cells.Item(4, 13).Value = -149.2  # ALC!M4 was -108
$ws.Cells.Item(4, 14).Value = -1075.25  # ALC!N4 was -1175

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 280.42856  # ALC!H33 was 318.18182
$ws.Cells.Item(33, 9).Value = 294.7  # ALC!I33 was 345
$ws.Cells.Item(33, 10).Value = 244.75  # ALC!J33 was 246.66667
$ws.Cells.Item(33, 11).Value = 294.7  # ALC!K33 was 345
$ws.Cells.Item(33, 12).Value = 244.75  # ALC!L33 was 246.66667
$ws.Cells.Item(33, 13).Value = -65.69999999999999  # ALC!M33 was -116
$ws.Cells.Item(33, 14).Value = -702.75  # ALC!N33 was -704.6666700000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 11335.333  # ALC!H69 was 10799.4
$ws.Cells.Item(69, 10).Value = 11335.333  # ALC!J69 was 10799.4
$ws.Cells.Item(69, 12).Value = 34005.999  # ALC!L69 was 32398.2
$ws.Cells.Item(69, 14).Value = -35753.999  # ALC!N69 was -34146.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 11335.333  # ALC!H72 was 10799.4
$ws.Cells.Item(72, 10).Value = 11335.333  # ALC!J72 was 10799.4
$ws.Cells.Item(72, 12).Value = 102017.997  # ALC!L72 was 97194.59999999999
$ws.Cells.Item(72, 14).Value = -110753.997  # ALC!N72 was -105930.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 7695052.5  # ALC!H76 was 9093817
$ws.Cells.Item(76, 9).Value = 11113865  # ALC!I76 was 12502873
$ws.Cells.Item(76, 10).Value = 2725.5  # ALC!J76 was 3000.6667
$ws.Cells.Item(76, 11).Value = 11113865  # ALC!K76 was 12502873
$ws.Cells.Item(76, 12).Value = 2725.5  # ALC!L76 was 3000.6667
$ws.Cells.Item(76, 13).Value = -11113550  # ALC!M76 was -12502558
$ws.Cells.Item(76, 14).Value = -3355.5  # ALC!N76 was -3630.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 7695052.5  # ALC!H79 was 9093817
$ws.Cells.Item(79, 9).Value = 11113865  # ALC!I79 was 12502873
$ws.Cells.Item(79, 10).Value = 2725.5  # ALC!J79 was 3000.6667
$ws.Cells.Item(79, 11).Value = 11113865  # ALC!K79 was 12502873
$ws.Cells.Item(79, 12).Value = 2725.5  # ALC!L79 was 3000.6667
$ws.Cells.Item(79, 13).Value = -11112773  # ALC!M79 was -12501781
$ws.Cells.Item(79, 14).Value = -4909.5  # ALC!N79 was -5184.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2999.8333  # ALC!H86 was 2999.8
$ws.Cells.Item(86, 10).Value = 2999.8  # ALC!J86 was 2999.75
$ws.Cells.Item(86, 12).Value = 2999.8  # ALC!L86 was 2999.75
$ws.Cells.Item(86, 14).Value = -5245.8  # ALC!N86 was -5245.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 2999.8333  # ALC!H89 was 2999.8
$ws.Cells.Item(89, 10).Value = 2999.8  # ALC!J89 was 2999.75
$ws.Cells.Item(89, 12).Value = 14999  # ALC!L89 was 14998.75
$ws.Cells.Item(89, 14).Value = -26231  # ALC!N89 was -26230.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 6029.1562  # ALC!H98 was 6204.6772
$ws.Cells.Item(98, 9).Value = 7653.625  # ALC!I98 was 7960.826
$ws.Cells.Item(98, 11).Value = 7653.625  # ALC!K98 was 7960.826
$ws.Cells.Item(98, 13).Value = -6155.625  # ALC!M98 was -6462.826

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 6029.1562  # ALC!H122 was 6204.6772
$ws.Cells.Item(122, 9).Value = 7653.625  # ALC!I122 was 7960.826
$ws.Cells.Item(122, 11).Value = 22960.875  # ALC!K122 was 23882.478
$ws.Cells.Item(122, 13).Value = -20510.875  # ALC!M122 was -21432.478

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 58332.13  # ARM!H32 was 70391.63
$ws.Cells.Item(32, 9).Value = 70508.414  # ARM!I32 was 91880.38
$ws.Cells.Item(32, 11).Value = 70508.414  # ARM!K32 was 91880.38
$ws.Cells.Item(32, 13).Value = -70221.414  # ARM!M32 was -91593.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2750  # ARM!H45 was 3100
$ws.Cells.Item(45, 9).Value = 1833.3334  # ARM!I45 was 2250
$ws.Cells.Item(45, 11).Value = 1833.3334  # ARM!K45 was 2250
$ws.Cells.Item(45, 13).Value = -1456.3334  # ARM!M45 was -1873

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 7861.9697  # ARM!H61 was 8255
$ws.Cells.Item(61, 9).Value = 6462.8  # ARM!I61 was 6984.222
$ws.Cells.Item(61, 11).Value = 6462.8  # ARM!K61 was 6984.222
$ws.Cells.Item(61, 13).Value = -6250.8  # ARM!M61 was -6772.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 9541.137000000001  # ARM!H63 was 9054.583000000001
$ws.Cells.Item(63, 9).Value = 0  # ARM!I63 was 3950
$ws.Cells.Item(63, 10).Value = 9541.137000000001  # ARM!J63 was 9518.637000000001
$ws.Cells.Item(63, 11).Value = 0  # ARM!K63 was 3950
$ws.Cells.Item(63, 12).Value = 9541.137000000001  # ARM!L63 was 9518.637000000001
$ws.Cells.Item(63, 13).ClearContents()  # ARM!M63 was -3264
$ws.Cells.Item(63, 14).Value = -10913.137  # ARM!N63 was -10890.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 9541.137000000001  # ARM!H66 was 9054.583000000001
$ws.Cells.Item(66, 9).Value = 0  # ARM!I66 was 3950
$ws.Cells.Item(66, 10).Value = 9541.137000000001  # ARM!J66 was 9518.637000000001
$ws.Cells.Item(66, 11).Value = 0  # ARM!K66 was 19750
$ws.Cells.Item(66, 12).Value = 47705.685  # ARM!L66 was 47593.185
$ws.Cells.Item(66, 13).ClearContents()  # ARM!M66 was -16318
$ws.Cells.Item(66, 14).Value = -54569.685  # ARM!N66 was -54457.185

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3019.8333  # ARM!H122 was 2639
$ws.Cells.Item(122, 9).Value = 1373  # ARM!I122 was 1412.1666
$ws.Cells.Item(122, 10).Value = 4666.6665  # ARM!J122 was 10000
$ws.Cells.Item(122, 11).Value = 4119  # ARM!K122 was 4236.4998
$ws.Cells.Item(122, 12).Value = 13999.9995  # ARM!L122 was 30000
$ws.Cells.Item(122, 13).Value = -1669  # ARM!M122 was -1786.4998
$ws.Cells.Item(122, 14).Value = -18899.9995  # ARM!N122 was -34900

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 8402.52  # ARM!H132 was 6807
$ws.Cells.Item(132, 9).Value = 6842.357  # ARM!I132 was 5016.409
$ws.Cells.Item(132, 11).Value = 20527.071  # ARM!K132 was 15049.227
$ws.Cells.Item(132, 13).Value = -17997.071  # ARM!M132 was -12519.227

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 7861.9697  # ARM!H136 was 8255
$ws.Cells.Item(136, 9).Value = 6462.8  # ARM!I136 was 6984.222
$ws.Cells.Item(136, 11).Value = 19388.4  # ARM!K136 was 20952.666
$ws.Cells.Item(136, 13).Value = -16838.4  # ARM!M136 was -18402.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2869.1738  # BSM!H20 was 2761.625
$ws.Cells.Item(20, 9).Value = 3240.4167  # BSM!I20 was 3036.923
$ws.Cells.Item(20, 10).Value = 2464.182  # BSM!J20 was 2436.2727
$ws.Cells.Item(20, 11).Value = 3240.4167  # BSM!K20 was 3036.923
$ws.Cells.Item(20, 12).Value = 2464.182  # BSM!L20 was 2436.2727
$ws.Cells.Item(20, 13).Value = -2993.4167  # BSM!M20 was -2789.923
$ws.Cells.Item(20, 14).Value = -2958.182  # BSM!N20 was -2930.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(109, 8).Value = 60000  # BSM!H109 was 0
$ws.Cells.Item(109, 10).Value = 60000  # BSM!J109 was 0
$ws.Cells.Item(109, 12).Value = 60000  # BSM!L109 was 0
$ws.Cells.Item(109, 14).Value = -62774  # BSM!N109 was None

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 5618.159  # BSM!H134 was 6136.7075
$ws.Cells.Item(134, 9).Value = 4125.697  # BSM!I134 was 4685.1333
$ws.Cells.Item(134, 11).Value = 12377.091  # BSM!K134 was 14055.3999
$ws.Cells.Item(134, 13).Value = -9842.091  # BSM!M134 was -11520.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 178.53847  # CRP!H7 was 159.55556
$ws.Cells.Item(7, 9).Value = 265.375  # CRP!I7 was 220.75
$ws.Cells.Item(7, 10).Value = 39.6  # CRP!J7 was 37.166668
$ws.Cells.Item(7, 11).Value = 265.375  # CRP!K7 was 220.75
$ws.Cells.Item(7, 12).Value = 39.6  # CRP!L7 was 37.166668
$ws.Cells.Item(7, 13).Value = -152.375  # CRP!M7 was -107.75
$ws.Cells.Item(7, 14).Value = -265.6  # CRP!N7 was -263.166668

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1913.2727  # CRP!H16 was 2116.6667
$ws.Cells.Item(16, 9).Value = 1906.25  # CRP!I16 was 2092.8572
$ws.Cells.Item(16, 10).Value = 1932  # CRP!J16 was 2200
$ws.Cells.Item(16, 11).Value = 1906.25  # CRP!K16 was 2092.8572
$ws.Cells.Item(16, 12).Value = 1932  # CRP!L16 was 2200
$ws.Cells.Item(16, 13).Value = -1619.25  # CRP!M16 was -1805.8572
$ws.Cells.Item(16, 14).Value = -2506  # CRP!N16 was -2774

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1913.2727  # CRP!H113 was 2116.6667
$ws.Cells.Item(113, 9).Value = 1906.25  # CRP!I113 was 2092.8572
$ws.Cells.Item(113, 10).Value = 1932  # CRP!J113 was 2200
$ws.Cells.Item(113, 11).Value = 1906.25  # CRP!K113 was 2092.8572
$ws.Cells.Item(113, 12).Value = 1932  # CRP!L113 was 2200
$ws.Cells.Item(113, 13).Value = 263.75  # CRP!M113 was 77.14280000000008
$ws.Cells.Item(113, 14).Value = -6272  # CRP!N113 was -6540

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 9750.429  # CRP!H134 was 13932.667
$ws.Cells.Item(134, 9).Value = 8887  # CRP!I134 was 14440.833
$ws.Cells.Item(134, 11).Value = 26661  # CRP!K134 was 43322.499
$ws.Cells.Item(134, 13).Value = -24126  # CRP!M134 was -40787.499

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 4001803  # CUL!H112 was 5001997
$ws.Cells.Item(112, 9).Value = 5001254  # CUL!I112 was 6667996
$ws.Cells.Item(112, 11).Value = 15003762  # CUL!K112 was 20003988
$ws.Cells.Item(112, 13).Value = -15002654  # CUL!M112 was -20002880

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(69, 8).Value = 26000  # GSM!H69 was 33425
$ws.Cells.Item(69, 10).Value = 26000  # GSM!J69 was 33425
$ws.Cells.Item(69, 12).Value = 26000  # GSM!L69 was 33425
$ws.Cells.Item(69, 14).Value = -27498  # GSM!N69 was -34923

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(72, 8).Value = 26000  # GSM!H72 was 33425
$ws.Cells.Item(72, 10).Value = 26000  # GSM!J72 was 33425
$ws.Cells.Item(72, 12).Value = 78000  # GSM!L72 was 100275
$ws.Cells.Item(72, 14).Value = -85488  # GSM!N72 was -107763

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(120, 8).Value = 60000  # GSM!H120 was 0
$ws.Cells.Item(120, 10).Value = 60000  # GSM!J120 was 0
$ws.Cells.Item(120, 12).Value = 60000  # GSM!L120 was 0
$ws.Cells.Item(120, 14).Value = -69676  # GSM!N120 was None

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(121, 8).Value = 49500  # GSM!H121 was 48966.668
$ws.Cells.Item(121, 10).Value = 49500  # GSM!J121 was 48966.668
$ws.Cells.Item(121, 12).Value = 49500  # GSM!L121 was 48966.668
$ws.Cells.Item(121, 14).Value = -52994  # GSM!N121 was -52460.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 1635  # LTW!H9 was 2386
$ws.Cells.Item(9, 9).Value = 724.75  # LTW!I9 was 933
$ws.Cells.Item(9, 10).Value = 3455.5  # LTW!J9 was 3475.75
$ws.Cells.Item(9, 11).Value = 724.75  # LTW!K9 was 933
$ws.Cells.Item(9, 12).Value = 3455.5  # LTW!L9 was 3475.75
$ws.Cells.Item(9, 13).Value = -500.75  # LTW!M9 was -709
$ws.Cells.Item(9, 14).Value = -3903.5  # LTW!N9 was -3923.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3883.0588  # LTW!H22 was 3912.5293
$ws.Cells.Item(22, 10).Value = 4505.591  # LTW!J22 was 4551.136
$ws.Cells.Item(22, 12).Value = 4505.591  # LTW!L22 was 4551.136
$ws.Cells.Item(22, 14).Value = -5095.591  # LTW!N22 was -5141.136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3883.0588  # LTW!H27 was 3912.5293
$ws.Cells.Item(27, 10).Value = 4505.591  # LTW!J27 was 4551.136
$ws.Cells.Item(27, 12).Value = 4505.591  # LTW!L27 was 4551.136
$ws.Cells.Item(27, 14).Value = -4719.591  # LTW!N27 was -4765.136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 6890.4287  # LTW!H46 was 7378.154
$ws.Cells.Item(46, 9).Value = 1220.2  # LTW!I46 was 1667
$ws.Cells.Item(46, 11).Value = 1220.2  # LTW!K46 was 1667
$ws.Cells.Item(46, 13).Value = -1032.2  # LTW!M46 was -1479

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 195.75  # LTW!H55 was 199.8421
$ws.Cells.Item(55, 9).Value = 137.125  # LTW!I55 was 138.625
$ws.Cells.Item(55, 10).Value = 234.83333  # LTW!J55 was 244.36363
$ws.Cells.Item(55, 11).Value = 137.125  # LTW!K55 was 138.625
$ws.Cells.Item(55, 12).Value = 234.83333  # LTW!L55 was 244.36363
$ws.Cells.Item(55, 13).Value = 35.875  # LTW!M55 was 34.375
$ws.Cells.Item(55, 14).Value = -580.8333299999999  # LTW!N55 was -590.3636300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3715.6155  # LTW!H61 was 3876.7693
$ws.Cells.Item(61, 9).Value = 3257  # LTW!I61 was 3556.2856
$ws.Cells.Item(61, 11).Value = 3257  # LTW!K61 was 3556.2856
$ws.Cells.Item(61, 13).Value = -3055  # LTW!M61 was -3354.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2756.6667  # LTW!H68 was 2846.4285
$ws.Cells.Item(68, 9).Value = 2205  # LTW!I68 was 2283.3333
$ws.Cells.Item(68, 11).Value = 2205  # LTW!K68 was 2283.3333
$ws.Cells.Item(68, 13).Value = -1456  # LTW!M68 was -1534.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 2756.6667  # LTW!H71 was 2846.4285
$ws.Cells.Item(71, 9).Value = 2205  # LTW!I71 was 2283.3333
$ws.Cells.Item(71, 11).Value = 11025  # LTW!K71 was 11416.6665
$ws.Cells.Item(71, 13).Value = -7281  # LTW!M71 was -7672.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3715.6155  # LTW!H113 was 3876.7693
$ws.Cells.Item(113, 9).Value = 3257  # LTW!I113 was 3556.2856
$ws.Cells.Item(113, 11).Value = 3257  # LTW!K113 was 3556.2856
$ws.Cells.Item(113, 13).Value = -1087  # LTW!M113 was -1386.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(124, 8).Value = 65557.25  # LTW!H124 was 61625.6
$ws.Cells.Item(124, 10).Value = 65557.25  # LTW!J124 was 61625.6
$ws.Cells.Item(124, 12).Value = 65557.25  # LTW!L124 was 61625.6
$ws.Cells.Item(124, 14).Value = -75377.25  # LTW!N124 was -71445.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(129, 8).Value = 78714.5  # LTW!H129 was 88000
$ws.Cells.Item(129, 10).Value = 78714.5  # LTW!J129 was 88000
$ws.Cells.Item(129, 12).Value = 78714.5  # LTW!L129 was 88000
$ws.Cells.Item(129, 14).Value = -88714.5  # LTW!N129 was -98000

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(131, 8).Value = 90784.8  # LTW!H131 was 100982.2
$ws.Cells.Item(131, 10).Value = 90784.8  # LTW!J131 was 100982.2
$ws.Cells.Item(131, 12).Value = 90784.8  # LTW!L131 was 100982.2
$ws.Cells.Item(131, 14).Value = -100864.8  # LTW!N131 was -111062.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 0  # WVR!H39 was 30000
$ws.Cells.Item(39, 10).Value = 0  # WVR!J39 was 30000
$ws.Cells.Item(39, 12).Value = 0  # WVR!L39 was 30000
$ws.Cells.Item(39, 14).ClearContents()  # WVR!N39 was -30826

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1784.2858  # WVR!H100 was 1826
$ws.Cells.Item(100, 9).Value = 1265.1666  # WVR!I100 was 1283.7059
$ws.Cells.Item(100, 11).Value = 2530.3332  # WVR!K100 was 2567.4118
$ws.Cells.Item(100, 13).Value = -1989.3332  # WVR!M100 was -2026.4118
